# Insert a new row above row 14 (shifts existing rows 14..118 down to 15..119)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the new weekly data point
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C14").Value = 'Arica y Parinacota'
$ws.Range("D14").Value = 45230
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = 100112012
$ws.Range("G14").Value = 'Espinaca'
$ws.Range("H14").Value = 'Sin especificar'
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 450
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = 1278
$ws.Range("N14").Value = '$/atado 2,5 a 3 kilos'
$ws.Range("O14").Value = 'Región de Arica y Parinacota'
$ws.Range("P14").Value = 426
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 'Hortaliza'
